# Update the build version timestamp throughout the workbook.
# Old timestamp: February 03 2026 17.29.55 EST
# New timestamp: February 03 2026 18.05.36 EST

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet updates ---

# A2: "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on ...)"
$a2 = $aboutSheet.Range("A2").Value2
$aboutSheet.Range("A2").Value2 = $a2.Replace($oldStamp, $newStamp)

# A6: "Recommended Citation: ... version '...(built on ...)' ..."
$a6 = $aboutSheet.Range("A6").Value2
$aboutSheet.Range("A6").Value2 = $a6.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet updates ---
# Column S ("build_version") for data rows 2 through 27.

for ($row = 2; $row -le 27; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    $val = $cell.Value2
    if ($val -ne $null -and $val.Contains($oldStamp)) {
        $cell.Value2 = $val.Replace($oldStamp, $newStamp)
    }
}
